$d = $word.ActiveDocument

# --- Edit 1: "We strive to..." paragraph -> "The self funding mechanism..." ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("We strive to", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find target paragraph 1" }
$p1 = $r1.Paragraphs(1)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/><w:b/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">The self funding mechanism that pays for development and encourages the community to decide how budget funds are spent.</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# --- Edit 2: "SmartCash mining prevents..." paragraph -> "SmartCash now has SmartMining..." ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("SmartCash mining prevents", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find target paragraph 2" }
$p2 = $r2.Paragraphs(1)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr><w:t xml:space="preserve">SmartCash now has SmartMining that prevents mining attacks. Mining can be done by anyone with a computer with one or more graphics cards. ASICs have yet to be created for the Keccak mining algorithm and it’s probably safe to assume no ASICs will be created for quite some time.</w:t></w:r></w:p>'
$p2.Range.InsertXML($xml2)
